# Auto-generated script to apply row-swap corrections to "Croatia 3NL" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: apply data previously belonging to row 9
$ws.Cells.Item(8, 2).Value2 = 7126860
$rowVals = @("NK Tomislav", "NK Oriolik Oriovac", 1, 1, "D", 2.5, 3.4, 2.4, 2.625, 3.4, 2.3, 0, 1.975, 1.725, 3, 1.95, 1.85, -1, 2.4, -1, 0, 0, -1, 0.8500000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F8:AC8").Value2 = $arr

# Row 9: apply data previously belonging to row 8
$ws.Cells.Item(9, 2).Value2 = 7126858
$rowVals = @("Sava Strmec", "Lucko", 0, 2, "A", 2.75, 3.4, 2.2, 3.75, 3.6, 1.8, 0.5, 1.95, 1.85, 2.75, 1.9, 1.9, -1, -1, 0.8, -1, 0.8500000000000001, -1, 0.8999999999999999)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F9:AC9").Value2 = $arr

# Row 14: apply data previously belonging to row 16
$ws.Cells.Item(14, 2).Value2 = 7183191
$rowVals = @("Zmaj Makarska", "NK Junak", 1, 2, "A", 3.6, 3.6, 1.8, 3.6, 3.6, 1.8, 0.5, 1.95, 1.85, 2.75, 1.9, 1.9, -1, -1, 0.8, -1, 0.8500000000000001, 0.45, -0.5)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F14:AC14").Value2 = $arr

# Row 16: apply data previously belonging to row 14
$ws.Cells.Item(16, 2).Value2 = 7183189
$rowVals = @("NK Vodice", "NK Uskok", 0, 1, "A", 2.1, 3.3, 3, 2.1, 3.3, 3, -0.25, 1.9, 1.9, 2.75, 1.925, 1.775, -1, -1, 2, -1, 0.8999999999999999, -1, 0.7749999999999999)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F16:AC16").Value2 = $arr

# Row 18: apply data previously belonging to row 19
$ws.Cells.Item(18, 2).Value2 = 7188989
$rowVals = @("Sava Strmec", "NK Vrapce", 1, 2, "A", 2.5, 3.4, 2.4, 2.5, 3.4, 2.4, 0, 1.95, 1.85, 3, 1.975, 1.825, -1, -1, 1.4, -1, 0.8500000000000001, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F18:AC18").Value2 = $arr

# Row 19: apply data previously belonging to row 18
$ws.Cells.Item(19, 2).Value2 = 7188990
$rowVals = @("NK Granicar Zupanja", "NK Bedem Ivankovo", 0, 2, "A", 2.1, 3.5, 3, 2.1, 3.5, 3, -0.25, 1.875, 1.925, 2.75, 1.825, 1.975, -1, -1, 2, -1, 0.925, -1, 0.9750000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F19:AC19").Value2 = $arr

# Row 22: apply data previously belonging to row 24
$ws.Cells.Item(22, 2).Value2 = 7202437
$rowVals = @("NK Maksimir", "Sava Strmec", 6, 0, "H", 1.4, 4.333, 6, 1.4, 4.333, 6, -1.25, 1.85, 1.95, 3, 1.8, 2, 0.3999999999999999, -1, -1, 0.8500000000000001, -1, 0.8, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F22:AC22").Value2 = $arr

# Row 23: apply data previously belonging to row 22
$ws.Cells.Item(23, 2).Value2 = 7202435
$rowVals = @("NK Udarnik Kurilovec", "NK Mladost Petrinja", 6, 0, "H", 2, 3.4, 3.1, 2, 3.4, 3.1, -0.25, 1.8, 2, 2.75, 1.825, 1.975, 1, -1, -1, 0.8, -1, 0.825, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F23:AC23").Value2 = $arr

# Row 24: apply data previously belonging to row 23
$ws.Cells.Item(24, 2).Value2 = 7202436
$rowVals = @("Lucko", "NK Tondach", 2, 1, "H", 1.615, 3.75, 4.333, 1.533, 4, 5, -0.75, 1.7, 2.1, 2.75, 1.825, 1.975, 0.5329999999999999, -1, -1, 0.35, -0.5, 0.4125, -0.5)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F24:AC24").Value2 = $arr

# Row 29: apply data previously belonging to row 30
$ws.Cells.Item(29, 2).Value2 = 7250138
$rowVals = @("NK Tomislav", "Sloga Nova Gradiska", 2, 2, "D", 2.1, 3.4, 2.9, 2.1, 3.4, 2.9, -0.25, 1.9, 1.9, 3, 1.9, 1.9, -1, 2.4, -1, -0.5, 0.45, 0.8999999999999999, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F29:AC29").Value2 = $arr

# Row 30: apply data previously belonging to row 29
$ws.Cells.Item(30, 2).Value2 = 7250137
$rowVals = @("NK Granicar Zupanja", "NK Svacic", 0, 2, "A", 1.727, 3.75, 3.75, 1.727, 3.75, 3.75, -0.75, 1.975, 1.825, 2.5, 1.8, 2, -1, -1, 2.75, -1, 0.825, -1, 1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F30:AC30").Value2 = $arr

# Row 33: apply data previously belonging to row 34
$ws.Cells.Item(33, 2).Value2 = 7291472
$rowVals = @("NK Lukavec", "Sava Strmec", 3, 0, "H", 2.2, 3.6, 2.6, 2.2, 3.6, 2.625, -0.25, 2, 1.8, 2.5, 1.8, 2, 1.2, -1, -1, 1, -1, 0.8, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F33:AC33").Value2 = $arr

# Row 34: apply data previously belonging to row 33
$ws.Cells.Item(34, 2).Value2 = 7291473
$rowVals = @("NK Maksimir", "NK Mladost Petrinja", 5, 0, "H", 1.25, 6, 7, 1.25, 6, 7.5, -1.75, 1.9, 1.9, 3, 1.975, 1.825, 0.25, -1, -1, 0.8999999999999999, -1, 0.9750000000000001, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F34:AC34").Value2 = $arr

# Row 35: apply data previously belonging to row 36
$ws.Cells.Item(35, 2).Value2 = 7305882
$rowVals = @("Zmaj Makarska", "NK Zadar", 0, 2, "A", 3.6, 4, 1.727, 3.6, 4, 1.727, 0.75, 1.825, 1.975, 2.75, 1.8, 2, -1, -1, 0.7270000000000001, -1, 0.9750000000000001, -1, 1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F35:AC35").Value2 = $arr

# Row 36: apply data previously belonging to row 35
$ws.Cells.Item(36, 2).Value2 = 7305881
$rowVals = @("NK Zagora", "NK Neretva", 4, 2, "H", 1.363, 4.75, 5.75, 1.363, 4.75, 5.75, -1.25, 1.75, 1.95, 3, 2, 1.8, 0.363, -1, -1, 0.75, -1, 1, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F36:AC36").Value2 = $arr

# Row 39: apply data previously belonging to row 40
$ws.Cells.Item(39, 2).Value2 = 7337011
$rowVals = @("NK Omis", "Zmaj Makarska", 3, 1, "H", 1.4, 4.75, 5.25, 1.4, 4.75, 5.25, -1.25, 1.9, 1.9, 2.75, 1.9, 1.9, 0.3999999999999999, -1, -1, 0.8999999999999999, -1, 0.8999999999999999, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F39:AC39").Value2 = $arr

# Row 40: apply data previously belonging to row 39
$ws.Cells.Item(40, 2).Value2 = 7337013
$rowVals = @("RNK Split", "NK Primorac Biograd", 1, 1, "D", 2.75, 3.4, 2.2, 2.75, 3.4, 2.2, 0.25, 1.8, 2, 3, 1.85, 1.95, -1, 2.4, -1, 0.4, -0.5, -1, 0.95)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F40:AC40").Value2 = $arr

# Row 41: apply data previously belonging to row 42
$ws.Cells.Item(41, 2).Value2 = 7337010
$rowVals = @("NK Maksimir", "NK Udarnik Kurilovec", 5, 0, "H", 1.5, 4, 5, 1.45, 4.2, 5.5, -1.25, 1.925, 1.875, 3, 1.75, 1.95, 0.45, -1, -1, 0.925, -1, 0.75, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F41:AC41").Value2 = $arr

# Row 42: apply data previously belonging to row 41
$ws.Cells.Item(42, 2).Value2 = 7337008
$rowVals = @("NK Dinamo Odranski Obre", "Gaj Mace", 0, 1, "A", 2.5, 3.4, 2.375, 3.2, 3.5, 1.95, 0.5, 1.75, 1.95, 3, 1.95, 1.85, -1, -1, 0.95, -1, 0.95, -1, 0.8500000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F42:AC42").Value2 = $arr

# Row 43: apply data previously belonging to row 44
$ws.Cells.Item(43, 2).Value2 = 7343496
$rowVals = @("NK Bedem Ivankovo", "Valpovka", 3, 0, "H", 1.85, 3.5, 3.5, 1.85, 3.5, 3.4, -0.5, 1.925, 1.875, 3, 1.875, 1.925, 0.8500000000000001, -1, -1, 0.925, -1, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F43:AC43").Value2 = $arr

# Row 44: apply data previously belonging to row 43
$ws.Cells.Item(44, 2).Value2 = 7343497
$rowVals = @("NK Tomislav", "NK Kutjevo", 3, 0, "H", 2.15, 3.6, 2.7, 2, 3.6, 3, -0.25, 1.8, 2, 3, 1.8, 2, 1, -1, -1, 0.8, -1, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F44:AC44").Value2 = $arr

# Row 50: apply data previously belonging to row 51
$ws.Cells.Item(50, 2).Value2 = 7382547
$rowVals = @("Zmaj Makarska", "RNK Split", 5, 0, "H", 2.3, 3.4, 2.625, 2.3, 3.4, 2.625, -0.25, 2.05, 1.75, 2.75, 1.975, 1.825, 1.3, -1, -1, 1.05, -1, 0.9750000000000001, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F50:AC50").Value2 = $arr

# Row 51: apply data previously belonging to row 50
$ws.Cells.Item(51, 2).Value2 = 7382546
$rowVals = @("NK Primorac Biograd", "NK Vodice", 1, 1, "D", 1.363, 4.333, 6.5, 1.25, 5, 9, -1.75, 1.975, 1.825, 2.75, 1.8, 2, -1, 4, -1, -1, 0.825, -1, 1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F51:AC51").Value2 = $arr

# Row 58: apply data previously belonging to row 59
$ws.Cells.Item(58, 2).Value2 = 7460852
$rowVals = @("NK Vodice", "NK HV Posedarje", 2, 2, "D", 4, 3.5, 1.727, 4, 3.5, 1.727, 0.75, 1.8, 2, 2.75, 1.95, 1.85, -1, 2.5, -1, 0.8, -1, 0.95, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F58:AC58").Value2 = $arr

# Row 59: apply data previously belonging to row 58
$ws.Cells.Item(59, 2).Value2 = 7460854
$rowVals = @("RNK Split", "NK Zadar", 1, 2, "A", 5.75, 4.5, 1.4, 5.5, 4.5, 1.4, 1.25, 1.95, 1.85, 3, 1.775, 2.025, -1, -1, 0.3999999999999999, 0.475, -0.5, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F59:AC59").Value2 = $arr

# Row 65: apply data previously belonging to row 67
$ws.Cells.Item(65, 2).Value2 = 7493772
$rowVals = @("NK Bedem Ivankovo", "Sloga Nova Gradiska", 2, 1, "H", 1.909, 3.6, 3.2, 1.4, 4.2, 6.5, -1.25, 1.9, 1.9, 3, 1.85, 1.95, 0.3999999999999999, -1, -1, -0.5, 0.45, 0, 0)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F65:AC65").Value2 = $arr

# Row 66: apply data previously belonging to row 68
$ws.Cells.Item(66, 2).Value2 = 7493775
$rowVals = @("NK Omis", "RNK Split", 4, 1, "H", 1.727, 3.75, 3.75, 1.5, 3.8, 5.25, -1, 1.925, 1.875, 2.5, 1.875, 1.925, 0.5, -1, -1, 0.925, -1, 0.875, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F66:AC66").Value2 = $arr

# Row 67: apply data previously belonging to row 65
$ws.Cells.Item(67, 2).Value2 = 7493774
$rowVals = @("NK Junak", "NK Gosk Kastel Gomilica", 1, 1, "D", 1.666, 3.75, 4, 1.45, 3.8, 5.75, -1, 1.825, 1.975, 2.75, 1.925, 1.875, -1, 2.8, -1, -1, 0.9750000000000001, -1, 0.875)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F67:AC67").Value2 = $arr

# Row 68: apply data previously belonging to row 66
$ws.Cells.Item(68, 2).Value2 = 7493773
$rowVals = @("NK Zadar", "NK Vodice", 2, 2, "D", 1.2, 6, 9, 1.2, 6, 9, -2, 1.85, 1.95, 3.75, 1.95, 1.85, -1, 5, -1, -1, 0.95, 0.475, -0.5)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F68:AC68").Value2 = $arr

# Row 73: apply data previously belonging to row 74
$ws.Cells.Item(73, 2).Value2 = 7517624
$rowVals = @("NK Croatia Dakovo", "Valpovka", 1, 6, "A", 1.4, 4.8, 5.25, 1.4, 4.8, 5.25, -1.25, 1.85, 1.95, 3.25, 1.85, 1.95, -1, -1, 4.25, -1, 0.95, 0.8500000000000001, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F73:AC73").Value2 = $arr

# Row 74: apply data previously belonging to row 73
$ws.Cells.Item(74, 2).Value2 = 7517623
$rowVals = @("NK Gosk Kastel Gomilica", "NK Primorac Biograd", 2, 0, "H", 2.05, 3.5, 3, 2.05, 3.5, 3, -0.25, 1.85, 1.95, 3, 1.975, 1.825, 1.05, -1, -1, 0.8500000000000001, -1, -1, 0.825)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F74:AC74").Value2 = $arr

# Row 78: apply data previously belonging to row 79
$ws.Cells.Item(78, 2).Value2 = 7519479
$rowVals = @("Sava Strmec", "NK Ponikve", 0, 1, "A", 2.75, 3.4, 2.2, 2.75, 3.4, 2.2, 0.25, 1.8, 2, 2.75, 1.9, 1.9, -1, -1, 1.2, -1, 1, -1, 0.8999999999999999)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F78:AC78").Value2 = $arr

# Row 79: apply data previously belonging to row 78
$ws.Cells.Item(79, 2).Value2 = 7519478
$rowVals = @("NK Bistra", "Lucko", 1, 2, "A", 3, 3.6, 2, 3, 3.6, 2, 0.25, 2, 1.8, 2.5, 1.8, 2, -1, -1, 1, -1, 0.8, 0.8, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F79:AC79").Value2 = $arr

# Row 89: apply data previously belonging to row 90
$ws.Cells.Item(89, 2).Value2 = 7939609
$rowVals = @("NK Neretvanac Opuzen", "NK GOSK Dubrovnik", 1, 0, "H", 1.4, 4.8, 5.25, 1.4, 5, 5, -1.25, 1.875, 1.925, 2.75, 1.8, 2, 0.3999999999999999, -1, -1, -0.5, 0.4625, -1, 1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F89:AC89").Value2 = $arr

# Row 90: apply data previously belonging to row 89
$ws.Cells.Item(90, 2).Value2 = 7939610
$rowVals = @("NK Sloga Mravince", "NK Zadar", 0, 4, "A", 2.75, 3.4, 2.2, 3.1, 3.4, 2, 0.25, 2, 1.8, 2.75, 1.8, 2, -1, -1, 1, -1, 0.8, 0.8, -1)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F90:AC90").Value2 = $arr

# Row 93: apply data previously belonging to row 94
$ws.Cells.Item(93, 2).Value2 = 7964849
$rowVals = @("NK Croatia Dakovo", "NK Tomislav", 2, 0, "H", 1.4, 4.5, 5.5, 1.4, 4.5, 5.5, -1.25, 1.85, 1.95, 3, 1.95, 1.85, 0.3999999999999999, -1, -1, 0.8500000000000001, -1, -1, 0.8500000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F93:AC93").Value2 = $arr

# Row 94: apply data previously belonging to row 93
$ws.Cells.Item(94, 2).Value2 = 7964850
$rowVals = @("NK Medulin", "NK Rudar Labin", 0, 0, "D", 3.8, 3.4, 1.8, 2.7, 3.4, 2.25, 0.25, 1.775, 2.025, 3, 1.95, 1.85, -1, 2.4, -1, 0.3875, -0.5, -1, 0.8500000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F94:AC94").Value2 = $arr

# Row 101: apply data previously belonging to row 102
$ws.Cells.Item(101, 2).Value2 = 8001390
$rowVals = @("Sava Strmec", "NK Maksimir", 1, 1, "D", 5.5, 4.5, 1.4, 3, 3.2, 2.15, 0.25, 1.875, 1.925, 2.5, 1.95, 1.85, -1, 2.2, -1, 0.4375, -0.5, -1, 0.8500000000000001)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F101:AC101").Value2 = $arr

# Row 102: apply data previously belonging to row 101
$ws.Cells.Item(102, 2).Value2 = 8001389
$rowVals = @("NK Bistra", "HNK Segesta", 1, 2, "A", 5.5, 4.5, 1.4, 3.8, 3.5, 1.75, 0.75, 1.75, 1.95, 2.75, 1.925, 1.775, -1, -1, 0.75, -0.5, 0.475, 0.4625, -0.5)
$arr = New-Object 'object[,]' 1,24
for ($i = 0; $i -lt 24; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("F102:AC102").Value2 = $arr
